# Apply cryptos list update (Fri Aug  2 05:37:18 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($ws, [string]$ref, [string]$val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "64.360.90"
$ws.Range("E2").Value = "  +0.34%  "
Set-TextValue $ws "D3" "3.159.78"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws "D5" "571.15"
$ws.Range("E5").Value = "  +0.31%  "
Set-TextValue $ws "D6" "163.76"
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("E7").Value = "  +0.07%  "
Set-TextValue $ws "D8" "0.577"
$ws.Range("E8").Value = "  -4.82%  "
$ws.Range("E9").Value = "  -2.83%  "
Set-TextValue $ws "D10" "6.61"
$ws.Range("E10").Value = "  -1.23%  "
Set-TextValue $ws "D11" "0.382"
$ws.Range("E11").Value = "  -0.46%  "
Set-TextValue $ws "D12" "3.718.43"
Set-TextValue $ws "D14" "64.431.52"
$ws.Range("E14").Value = "  +0.31%  "
Set-TextValue $ws "D15" "25.17"
$ws.Range("E15").Value = "  -0.64%  "
Set-TextValue $ws "D16" "3.157.18"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("E17").Value = "  -2.29%  "
Set-TextValue $ws "D18" "404.75"
$ws.Range("E18").Value = "  -2.82%  "
Set-TextValue $ws "D19" "12.68"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("E21").Value = "  +0.52%  "
Set-TextValue $ws "D23" "68.56"
$ws.Range("E23").Value = "  -2.08%  "
Set-TextValue $ws "D24" "0.484"
$ws.Range("E24").Value = "  -1.30%  "
Set-TextValue $ws "D25" "0.193"
$ws.Range("E25").Value = "  -4.40%  "
$ws.Range("E26").Value = "  -3.98%  "
Set-TextValue $ws "D27" "8.82"
$ws.Range("E27").Value = "  +0.70%  "
$ws.Range("E28").Value = "  +0.05%  "
Set-TextValue $ws "D29" "1.82"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -2.66%  "
Set-TextValue $ws "D31" "6.29"
$ws.Range("E31").Value = "  -0.76%  "
Set-TextValue $ws "D32" "4.84"
$ws.Range("E32").Value = "  -3.43%  "
Set-TextValue $ws "D33" "156.54"
$ws.Range("E33").Value = "  +0.85%  "
Set-TextValue $ws "D34" "1.11"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  -2.57%  "
Set-TextValue $ws "D36" "2.674.90"
$ws.Range("E36").Value = "  -0.87%  "
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("E38").Value = "  -2.78%  "
Set-TextValue $ws "D39" "4.08"
$ws.Range("E39").Value = "  -2.03%  "
$ws.Range("E40").Value = "  -2.00%  "
Set-TextValue $ws "D41" "0.0616"
$ws.Range("E41").Value = "  -1.19%  "
Set-TextValue $ws "D42" "5.44"
$ws.Range("E42").Value = "  -3.62%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws "D43" "0.0256"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D44" "288.88"
$ws.Range("E44").Value = "  -2.17%  "
Set-TextValue $ws "D45" "21.25"
$ws.Range("E45").Value = "  -2.56%  "
Set-TextValue $ws "D46" "1.00"
$ws.Range("E46").Value = "  +0.02%  "
Set-TextValue $ws "D47" "0.0982"
$ws.Range("E47").Value = "  -1.00%  "
Set-TextValue $ws "D48" "10.52"
$ws.Range("E48").Value = "  +0.72%  "
$ws.Range("E49").Value = "  -7.51%  "
Set-TextValue $ws "D50" "5.69"
$ws.Range("E50").Value = "  -1.46%  "
Set-TextValue $ws "D51" "0.877"
$ws.Range("E51").Value = "  -6.34%  "
